# Weekly update: insert two new "Cebollín" price rows (dated 2023-06-16 / serial 45093)
# right above what used to be row 123, pushing the existing rows 123-157 down to 125-159.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 123-124; everything that used to start at row 123
# (through row 157) shifts down to rows 125-159, and the sheet dimension grows
# from A1:R157 to A1:R159 automatically.
$ws.Rows("123:124").Insert()

# Common / constant values shared by all data rows on this sheet.
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$catId     = 100112037
$categoria = "Cebollín"
$variedad  = "Sin especificar"
$origen    = "Provincia de Diguillín"
$clasif    = "Hortaliza"

# New row 123 - "Primera" quality
$ws.Cells.Item(123, 1).Value  = $mercadoId
$ws.Cells.Item(123, 2).Value  = $mercado
$ws.Cells.Item(123, 3).Value  = $region
$ws.Cells.Item(123, 4).Value  = 45093
$ws.Cells.Item(123, 5).Value  = $codreg
$ws.Cells.Item(123, 6).Value  = $catId
$ws.Cells.Item(123, 7).Value  = $categoria
$ws.Cells.Item(123, 8).Value  = $variedad
$ws.Cells.Item(123, 9).Value  = "Primera"
$ws.Cells.Item(123, 10).Value = 60
$ws.Cells.Item(123, 11).Value = 6000
$ws.Cells.Item(123, 12).Value = 6000
$ws.Cells.Item(123, 13).Value = 6000
$ws.Cells.Item(123, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(123, 15).Value = $origen
$ws.Cells.Item(123, 16).Value = 167
$ws.Cells.Item(123, 17).Value = 36
$ws.Cells.Item(123, 18).Value = $clasif

# New row 124 - "Segunda" quality
$ws.Cells.Item(124, 1).Value  = $mercadoId
$ws.Cells.Item(124, 2).Value  = $mercado
$ws.Cells.Item(124, 3).Value  = $region
$ws.Cells.Item(124, 4).Value  = 45093
$ws.Cells.Item(124, 5).Value  = $codreg
$ws.Cells.Item(124, 6).Value  = $catId
$ws.Cells.Item(124, 7).Value  = $categoria
$ws.Cells.Item(124, 8).Value  = $variedad
$ws.Cells.Item(124, 9).Value  = "Segunda"
$ws.Cells.Item(124, 10).Value = 60
$ws.Cells.Item(124, 11).Value = 5000
$ws.Cells.Item(124, 12).Value = 5000
$ws.Cells.Item(124, 13).Value = 5000
$ws.Cells.Item(124, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(124, 15).Value = $origen
$ws.Cells.Item(124, 16).Value = 139
$ws.Cells.Item(124, 17).Value = 36
$ws.Cells.Item(124, 18).Value = $clasif
